$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F ("想去人数")
$updates = @{
    2  = 853
    4  = 2165
    6  = 12591
    7  = 58
    9  = 504
    10 = 460
    11 = 1149
    12 = 947
    13 = 13657
    14 = 13973
    19 = 8
    23 = 1052
    26 = 585
    27 = 5088
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
